# Scheduled market-data refresh: update computed Leve price/profit columns (H:N)
# across all job sheets in the Adamantoise_Profits workbook.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2510.25
$ws.Range("I6").Value = 2510.25
$ws.Range("K6").Value = 7530.75
$ws.Range("M6").Value = -7418.75
$ws.Range("H12").Value = 28476.857
$ws.Range("I12").Value = 188.66667
$ws.Range("J12").Value = 49693
$ws.Range("K12").Value = 188.66667
$ws.Range("L12").Value = 49693
$ws.Range("M12").Value = -18.66667000000001
$ws.Range("N12").Value = -50033
$ws.Range("H15").Value = 1154.3768
$ws.Range("I15").Value = 1154.3768
$ws.Range("K15").Value = 3463.1304
$ws.Range("M15").Value = -3294.1304
$ws.Range("H38").Value = 2270.25
$ws.Range("I38").Value = 193.66667
$ws.Range("J38").Value = 8500
$ws.Range("K38").Value = 581.00001
$ws.Range("L38").Value = 25500
$ws.Range("M38").Value = -209.00001
$ws.Range("N38").Value = -26244
$ws.Range("H41").Value = 1842.9
$ws.Range("I41").Value = 2072.7144
$ws.Range("J41").Value = 1306.6666
$ws.Range("K41").Value = 2072.7144
$ws.Range("L41").Value = 1306.6666
$ws.Range("M41").Value = -1632.7144
$ws.Range("N41").Value = -2186.6666
$ws.Range("H55").Value = 77052.84
$ws.Range("I55").Value = 50
$ws.Range("J55").Value = 111276.336
$ws.Range("K55").Value = 50
$ws.Range("L55").Value = 111276.336
$ws.Range("M55").Value = 164
$ws.Range("N55").Value = -111704.336
$ws.Range("H86").Value = 166669220
$ws.Range("I86").Value = 125002090
$ws.Range("K86").Value = 125002090
$ws.Range("M86").Value = -125000967
$ws.Range("H89").Value = 166669220
$ws.Range("I89").Value = 125002090
$ws.Range("K89").Value = 625010450
$ws.Range("M89").Value = -625004834
$ws.Range("H117").Value = 107967.25
$ws.Range("J117").Value = 107967.25
$ws.Range("L117").Value = 107967.25
$ws.Range("N117").Value = -117145.25
$ws.Range("H127").Value = 1687.9
$ws.Range("I127").Value = 1096.5
$ws.Range("K127").Value = 3289.5
$ws.Range("M127").Value = 1670.5
$ws.Range("H132").Value = 5512.5557
$ws.Range("I132").Value = 5826.625
$ws.Range("K132").Value = 17479.875
$ws.Range("M132").Value = -14949.875
$ws.Range("H137").Value = 166538
$ws.Range("I137").Value = 275533.34
$ws.Range("J137").Value = 3045
$ws.Range("K137").Value = 826600.02
$ws.Range("L137").Value = 9135
$ws.Range("M137").Value = -824050.02
$ws.Range("N137").Value = -14235

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2130.5881
$ws.Range("I2").Value = 1393.3334
$ws.Range("J2").Value = 3900
$ws.Range("K2").Value = 1393.3334
$ws.Range("L2").Value = 3900
$ws.Range("M2").Value = -1280.3334
$ws.Range("N2").Value = -4126
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()
$ws.Range("H116").Value = 2130.5881
$ws.Range("I116").Value = 1393.3334
$ws.Range("J116").Value = 3900
$ws.Range("K116").Value = 1393.3334
$ws.Range("L116").Value = 3900
$ws.Range("M116").Value = 900.6666
$ws.Range("N116").Value = -8488
$ws.Range("H122").Value = 4154.2896
$ws.Range("I122").Value = 2622.84
$ws.Range("K122").Value = 7868.52
$ws.Range("M122").Value = -5418.52

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2130.5881
$ws.Range("I3").Value = 1393.3334
$ws.Range("J3").Value = 3900
$ws.Range("K3").Value = 1393.3334
$ws.Range("L3").Value = 3900
$ws.Range("M3").Value = -1279.3334
$ws.Range("N3").Value = -4128
$ws.Range("H20").Value = 2822.3076
$ws.Range("I20").Value = 2383.7144
$ws.Range("K20").Value = 2383.7144
$ws.Range("M20").Value = -2136.7144

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3352.4473
$ws.Range("J31").Value = 5152.0713
$ws.Range("L31").Value = 5152.0713
$ws.Range("N31").Value = -5742.0713
$ws.Range("H34").Value = 3352.4473
$ws.Range("J34").Value = 5152.0713
$ws.Range("L34").Value = 5152.0713
$ws.Range("N34").Value = -5556.0713
$ws.Range("H107").Value = 1996
$ws.Range("J107").Value = 3737.5
$ws.Range("L107").Value = 3737.5
$ws.Range("N107").Value = -7577.5
$ws.Range("H122").Value = 1360.8182
$ws.Range("I122").Value = 1418.7778
$ws.Range("K122").Value = 4256.3334
$ws.Range("M122").Value = -1806.3334
$ws.Range("H132").Value = 1331.3334
$ws.Range("I132").Value = 1443.4615
$ws.Range("K132").Value = 4330.3845
$ws.Range("M132").Value = -1800.3845

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6472.25
$ws.Range("I3").Value = 6472.25
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 19416.75
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -19304.75
$ws.Range("N3").ClearContents()
$ws.Range("H9").Value = 3000
$ws.Range("J9").Value = 4000
$ws.Range("L9").Value = 12000
$ws.Range("N9").Value = -12448
$ws.Range("H38").Value = 136.77777
$ws.Range("J38").Value = 267.66666
$ws.Range("L38").Value = 802.9999799999999
$ws.Range("N38").Value = -1496.99998
$ws.Range("H50").Value = 2049.889
$ws.Range("I50").Value = 1170
$ws.Range("K50").Value = 3510
$ws.Range("M50").Value = -3029
$ws.Range("H53").Value = 2049.889
$ws.Range("I53").Value = 1170
$ws.Range("K53").Value = 3510
$ws.Range("M53").Value = -3029
$ws.Range("H75").Value = 749.5
$ws.Range("J75").Value = 500
$ws.Range("L75").Value = 1500
$ws.Range("N75").Value = -3496
$ws.Range("H78").Value = 749.5
$ws.Range("J78").Value = 500
$ws.Range("L78").Value = 4500
$ws.Range("N78").Value = -14484
$ws.Range("H133").Value = 3857.5
$ws.Range("I133").Value = 4476.6665
$ws.Range("K133").Value = 13429.9995
$ws.Range("M133").Value = -8369.999500000002
$ws.Range("H139").Value = 2305.7144
$ws.Range("I139").Value = 1990.6666
$ws.Range("K139").Value = 5971.9998
$ws.Range("M139").Value = -831.9997999999996

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4511.881
$ws.Range("J70").Value = 4499.974
$ws.Range("L70").Value = 4499.974
$ws.Range("N70").Value = -5039.974
$ws.Range("H73").Value = 4511.881
$ws.Range("J73").Value = 4499.974
$ws.Range("L73").Value = 4499.974
$ws.Range("N73").Value = -6371.974
$ws.Range("H122").Value = 1865.375
$ws.Range("I122").Value = 1641.1666
$ws.Range("K122").Value = 4923.4998
$ws.Range("M122").Value = -2473.4998

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 254.5
$ws.Range("I55").Value = 229.1875
$ws.Range("K55").Value = 229.1875
$ws.Range("M55").Value = -56.1875
$ws.Range("H99").Value = 78393.8
$ws.Range("I99").Value = 45000
$ws.Range("J99").Value = 100656.336
$ws.Range("K99").Value = 45000
$ws.Range("L99").Value = 100656.336
$ws.Range("M99").Value = -42005
$ws.Range("N99").Value = -106646.336
$ws.Range("H122").Value = 24019.23
$ws.Range("I122").Value = 22152.105
$ws.Range("K122").Value = 66456.315
$ws.Range("M122").Value = -64006.315
$ws.Range("H141").Value = 524999.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 524999.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 524999.5
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -535359.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 12678.1875
$ws.Range("I96").Value = 8012.25
$ws.Range("J96").Value = 17344.125
$ws.Range("K96").Value = 8012.25
$ws.Range("L96").Value = 17344.125
$ws.Range("M96").Value = -6639.25
$ws.Range("N96").Value = -20090.125
$ws.Range("H101").Value = 68425.25
$ws.Range("J101").Value = 68425.25
$ws.Range("L101").Value = 68425.25
$ws.Range("N101").Value = -74915.25
$ws.Range("H113").Value = 333.26666
$ws.Range("I113").Value = 364.36365
$ws.Range("K113").Value = 1093.09095
$ws.Range("M113").Value = 1076.90905
$ws.Range("H122").Value = 3688.4443
$ws.Range("I122").Value = 1623.9445
$ws.Range("J122").Value = 7817.4443
$ws.Range("K122").Value = 4871.833500000001
$ws.Range("L122").Value = 23452.3329
$ws.Range("M122").Value = -2421.833500000001
$ws.Range("N122").Value = -28352.3329
